# Generate Report for Handback
# For the localization rows (row 2: 1500ab6d..., row 3: 621bf619...) on the
# zh-cn and de-de sheets, record that the handback has happened:
#   - Status column (B) flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" (this text lives in the shared
#     string table and is also referenced by the Overview sheet, so that
#     sheet's Status cells update automatically).
#   - Latest Target File (E) / Latest Handback File (F) are now populated
#     with the same file names as the handoff columns (A / C).
#   - Latest Handback DateTime (G) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# The Overview sheet's Status columns (B, C) reuse the very same shared
# string as the per-language sheets' Status column (B), so flip them too
# -- that keeps everything pointed at one shared-string entry, exactly
# like the source edit did.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 2).Value = $newStatus
$wsOverview.Cells.Item(2, 3).Value = $newStatus
$wsOverview.Cells.Item(3, 2).Value = $newStatus
$wsOverview.Cells.Item(3, 3).Value = $newStatus

$sheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-02-17 04:57:51" },
    @{ Name = "de-de"; HandbackTime = "2016-02-17 04:58:10" }
)

foreach ($entry in $sheets) {
    $ws = $wb.Worksheets.Item($entry.Name)

    foreach ($row in 2, 3) {
        # Status -> Handed back
        $ws.Cells.Item($row, 2).Value = $newStatus

        # Latest Target File (E) mirrors Source File Name (A)
        $sourceFile = $ws.Cells.Item($row, 1).Text
        $ws.Cells.Item($row, 5).Value = $sourceFile

        # Latest Handback File (F) mirrors Latest Handoff File (C)
        $handoffFile = $ws.Cells.Item($row, 3).Text
        $ws.Cells.Item($row, 6).Value = $handoffFile

        # Apply the hyperlink style (same as used on A/C) to the new cells
        $ws.Range($ws.Cells.Item($row, 5), $ws.Cells.Item($row, 6)).Style = "HyperLink"

        # Hyperlinks for the newly-filled cells point at the same targets
        # as the corresponding handoff-file hyperlinks.
        $srcLink = $ws.Hyperlinks.Item($ws.Cells.Item($row, 1).Address())
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $srcLink.Address, "", "", $sourceFile)

        $handoffLink = $ws.Hyperlinks.Item($ws.Cells.Item($row, 3).Address())
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $handoffLink.Address, "", "", $handoffFile)

        # Latest Handback DateTime (G)
        $ws.Cells.Item($row, 7).Value = $entry.HandbackTime
    }
}
